$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date-column formatting from the last existing data row (343) down through the new rows (344-357)
# Only column A carries the special (bordered/centered/date-formatted) style; B/C/D stay default.
$ws.Range("A343").Copy() | Out-Null
$ws.Range("A344:A357").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(344, 1).Value = 44418
$ws.Cells.Item(344, 2).Value = 1
$ws.Cells.Item(344, 3).Value = 11
$ws.Cells.Item(344, 4).Value = 72.91045270762909

$ws.Cells.Item(345, 1).Value = 44419
$ws.Cells.Item(345, 2).Value = 3
$ws.Cells.Item(345, 3).Value = 14
$ws.Cells.Item(345, 4).Value = 92.79512162789156

$ws.Cells.Item(346, 1).Value = 44420
$ws.Cells.Item(346, 2).Value = 6
$ws.Cells.Item(346, 3).Value = 18
$ws.Cells.Item(346, 4).Value = 119.3080135215749

$ws.Cells.Item(347, 1).Value = 44421
$ws.Cells.Item(347, 2).Value = 5
$ws.Cells.Item(347, 3).Value = 20
$ws.Cells.Item(347, 4).Value = 132.5644594684165

$ws.Cells.Item(348, 1).Value = 44422
$ws.Cells.Item(348, 2).Value = 2
$ws.Cells.Item(348, 3).Value = 21
$ws.Cells.Item(348, 4).Value = 139.1926824418373

$ws.Cells.Item(349, 1).Value = 44423
$ws.Cells.Item(349, 2).Value = 0
$ws.Cells.Item(349, 3).Value = 20
$ws.Cells.Item(349, 4).Value = 132.5644594684165

$ws.Cells.Item(350, 1).Value = 44424
$ws.Cells.Item(350, 2).Value = 1
$ws.Cells.Item(350, 3).Value = 18
$ws.Cells.Item(350, 4).Value = 119.3080135215749

$ws.Cells.Item(351, 1).Value = 44425
$ws.Cells.Item(351, 2).Value = 5
$ws.Cells.Item(351, 3).Value = 22
$ws.Cells.Item(351, 4).Value = 145.8209054152582

$ws.Cells.Item(352, 1).Value = 44426
$ws.Cells.Item(352, 2).Value = 1
$ws.Cells.Item(352, 3).Value = 20
$ws.Cells.Item(352, 4).Value = 132.5644594684165

$ws.Cells.Item(353, 1).Value = 44427
$ws.Cells.Item(353, 2).Value = 1
$ws.Cells.Item(353, 3).Value = 15
$ws.Cells.Item(353, 4).Value = 99.42334460131239

$ws.Cells.Item(354, 1).Value = 44428
$ws.Cells.Item(354, 2).Value = 2
$ws.Cells.Item(354, 3).Value = 12
$ws.Cells.Item(354, 4).Value = 79.53867568104991

$ws.Cells.Item(355, 1).Value = 44429
$ws.Cells.Item(355, 2).Value = 1
$ws.Cells.Item(355, 3).Value = 11
$ws.Cells.Item(355, 4).Value = 72.91045270762909

$ws.Cells.Item(356, 1).Value = 44430
$ws.Cells.Item(356, 2).Value = 1
$ws.Cells.Item(356, 3).Value = 12
$ws.Cells.Item(356, 4).Value = 79.53867568104991

$ws.Cells.Item(357, 1).Value = 44431
$ws.Cells.Item(357, 2).Value = 0
$ws.Cells.Item(357, 3).Value = 11
$ws.Cells.Item(357, 4).Value = 72.91045270762909
